$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72; existing rows 72-161 shift down to 73-162.
$ws.Rows(72).Insert()

# Populate the newly inserted row 72 with the new data record.
$ws.Range("A72").Value = 8
$ws.Range("B72").Value = "Terminal La Palmera de La Serena"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44413
$ws.Range("E72").Value = 4
$ws.Range("F72").Value = 100114013
$ws.Range("G72").Value = "Zanahoria"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 720
$ws.Range("K72").Value = 5000
$ws.Range("L72").Value = 5500
$ws.Range("M72").Value = 5250
$ws.Range("N72").Value = "$/saco 20 kilos"
$ws.Range("O72").Value = "Provincia del Elquí"
$ws.Range("P72").Value = 262
$ws.Range("Q72").Value = 20
$ws.Range("R72").Value = "Hortaliza"
